$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:ALMONDZ"
$ws.Range("C2").Value = "NSE:ALKYLAMINE"
$ws.Range("E2").Value = "NSE:HDFCBANK"
$ws.Range("F2").Value = "NSE:ASIANPAINT"
$ws.Range("B3").Value = "NSE:APOLLOHOSP"
$ws.Range("C3").Value = "NSE:AMIORG"
$ws.Range("E3").Value = "NSE:LTTS"
$ws.Range("F3").Value = "NSE:GLENMARK"
$ws.Range("B4").Value = "NSE:AYMSYNTEX"
$ws.Range("C4").Value = "NSE:AWL"
$ws.Range("F4").Value = "NSE:IOC"
$ws.Range("B5").Value = "NSE:GLENMARK"
$ws.Range("C5").Value = "NSE:CAPTRUST"
$ws.Range("B6").Value = "NSE:GSLSU"
$ws.Range("C6").Value = "NSE:CHALET"
$ws.Range("B7").Value = "NSE:HINDPETRO"
$ws.Range("C7").Value = "NSE:CHAMBLFERT"
$ws.Range("B8").Value = "NSE:INDIACEM"
$ws.Range("C8").Value = "NSE:DCAL"
$ws.Range("B9").Value = "NSE:IZMO"
$ws.Range("C9").Value = "NSE:GTLINFRA"
$ws.Range("B10").Value = "NSE:KIRLPNU"
$ws.Range("C10").Value = "NSE:INDIANCARD"
$ws.Range("B11").Value = "NSE:MAANALU"
$ws.Range("C11").Value = "NSE:IRISDOREME"
$ws.Range("B12").Value = "NSE:NHPC"
$ws.Range("C12").Value = "NSE:KAMATHOTEL"
$ws.Range("B13").Value = "NSE:OIL"
$ws.Range("C13").Value = "NSE:KOPRAN"
$ws.Range("B14").Value = "NSE:REFEX"
$ws.Range("C14").Value = "NSE:LEXUS"
$ws.Range("B15").Value = "NSE:SAKSOFT"
$ws.Range("C15").Value = "NSE:MBLINFRA"
$ws.Range("C16").Value = "NSE:NGLFINE"
$ws.Range("C17").Value = "NSE:NINSYS"
$ws.Range("C18").Value = "NSE:OLECTRA"
$ws.Range("C19").Value = "NSE:RAMANEWS"
$ws.Range("C20").Value = "NSE:RKFORGE"

$ws.Range("A21:A33").EntireRow.Delete()
